$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name labels (column A) whose sort position shifted ---
$ws.Range("A7").Value = "Francia"
$ws.Range("A8").Value = "Alemania"
$ws.Range("A39").Value = "Peru"
$ws.Range("A40").Value = "Arabia Saudita"
$ws.Range("A41").Value = "Indonesia"
$ws.Range("A135").Value = "Aruba"
$ws.Range("A136").Value = "Guayana Francesa"

# --- Update numeric statistics (columns B-H) ---
$ws.Range("B4").Value = 386104
$ws.Range("C4").Value = 19100
$ws.Range("D4").Value = 21316
$ws.Range("E4").Value = 352546
$ws.Range("F4").Value = 9104
$ws.Range("G4").Value = 1371
$ws.Range("H4").Value = 12242

$ws.Range("B7").Value = 109069
$ws.Range("C7").Value = 11059
$ws.Range("D7").Value = 19337
$ws.Range("E7").Value = 79404
$ws.Range("F7").Value = 7131
$ws.Range("G7").Value = 1417
$ws.Range("H7").Value = 10328

$ws.Range("B8").Value = 106504
$ws.Range("C8").Value = 3129
$ws.Range("D8").Value = 36081
$ws.Range("E8").Value = 68482
$ws.Range("F8").Value = 4895
$ws.Range("G8").Value = 131
$ws.Range("H8").Value = 1941

$ws.Range("B16").Value = 17825
$ws.Range("C16").Value = 1158
$ws.Range("D16").Value = 3922
$ws.Range("E16").Value = 13529
$ws.Range("F16").Value = 426
$ws.Range("G16").Value = 51
$ws.Range("H16").Value = 374

$ws.Range("B25").Value = 5903
$ws.Range("C25").Value = 38
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = 5782
$ws.Range("F25").Value = 83

$ws.Range("B39").Value = 2954
$ws.Range("C39").Value = 393
$ws.Range("D39").Value = 997
$ws.Range("E39").Value = 1865
$ws.Range("F39").Value = 89
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 92

$ws.Range("B40").Value = 2795
$ws.Range("C40").Value = 190
$ws.Range("D40").Value = 615
$ws.Range("E40").Value = 2139
$ws.Range("F40").Value = 41
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 41

$ws.Range("B41").Value = 2738
$ws.Range("C41").Value = 247
$ws.Range("D41").Value = 204
$ws.Range("E41").Value = 2313
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 12
$ws.Range("H41").Value = 221

$ws.Range("B135").Value = 74
$ws.Range("C135").Value = 3
$ws.Range("D135").Value = 14
$ws.Range("E135").Value = 60
$ws.Range("F135").Value = 0

$ws.Range("B136").Value = 72
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 34
$ws.Range("E136").Value = 38
$ws.Range("F136").Value = 1

$ws.Range("B138").Value = 63
$ws.Range("C138").Value = 3
$ws.Range("D138").Value = 6
$ws.Range("E138").Value = 54
$ws.Range("F138").Value = 4

$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 3
$ws.Range("E163").Value = 14
$ws.Range("F163").Value = 0

# --- Update "last updated" timestamp message (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 19:52"
